$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "erros"

$ws.Range("A1").Value = "id"
$ws.Range("A2").Value = "DEV01"
$ws.Range("B2").Value = "Saldo não é igual a zero: xxx,xx débito: xxx,xx Crédito: xxx,xx"
$ws.Range("B1").Value = "erro"

$ws.Columns.Item(2).AutoFit() | Out-Null

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
$ws.PageSetup.LeftMargin = 36.850393728
$ws.PageSetup.RightMargin = 36.850393728
$ws.PageSetup.TopMargin = 56.692913399999995
$ws.PageSetup.BottomMargin = 56.692913399999995
$ws.PageSetup.HeaderMargin = 22.67716464
$ws.PageSetup.FooterMargin = 22.67716464

$ws.Range("B2").Select() | Out-Null
